# Applies the authored edit: the "Five Card Draw" title run gets merged
# back into a single run, the game-description slide gets a couple of
# small typo fixes, and the URL mentioned on the credits slide becomes a
# clickable hyperlink ("link works as one now").

$p = $ppt.ActivePresentation

# --- Slide 1: "Five Card " + "Draw" were two separate runs; PowerPoint's
#     spell-checker re-consolidated them into a single run once the file
#     was re-saved. Re-set the whole range so the engine rebuilds one run.
$s1  = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(2)
$tr1 = $sh1.TextFrame.TextRange
$full1 = $tr1.Text
$tr1.Characters(1, $full1.Length).Text = $full1

# --- Slide 2: fix two small typos in the game description ---
$s2  = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange

$full2 = $tr2.Text
$needle = "apoker"
$idx = $full2.IndexOf($needle)
$tr2.Characters($idx + 1, $needle.Length).Text = "poker"

$full2 = $tr2.Text
$needle = "draw again (with new cards)."
$idx = $full2.IndexOf($needle)
$tr2.Characters($idx + 1, $needle.Length).Text = "redraw again (with new cards from the deck)."

# --- Slide 13: turn the plain URL text into a real hyperlink ---
$s13  = $p.Slides.Item(13)
$sh13 = $s13.Shapes.Item(3)
$tr13 = $sh13.TextFrame.TextRange

$url = "https://github.com/ArchDevilCSharp2/Team-Work---Console-Game-Arch-Devil"
$full13 = $tr13.Text
$idx = $full13.IndexOf($url)
$urlRange = $tr13.Characters($idx + 1, $url.Length)
$urlRange.ActionSettings.Item(1).Hyperlink.Address = $url
